$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = "2025/12/04 08:00"
$ws.Range("B40").Value = "32,744位本"
$ws.Range("C40").Value = "87位 広告・宣伝 (本)"
$ws.Range("D40").Value = "140位商業デザイン"
$ws.Range("E40").Value = "1,749位ビジネス実用本"
$ws.Range("F40").Value = "-"
$ws.Range("G40").Value = "-"
